$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Text columns (A-D): temporarily force a text number format so the
# date-like / numeric-like strings are stored as literal text (matching
# how the existing rows above store Date/Time/Weekday/Week as text),
# instead of being auto-parsed into a date serial number or a plain
# number. Then restore the default formatting on those cells.
$textRange = $ws.Range("A9:D9")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-05-30"
$ws.Cells.Item($row, 2).Value = "22:07:01"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "22"

$textRange.ClearFormats()

# Numeric columns (E-T)
$ws.Cells.Item($row, 5).Value = 119934
$ws.Cells.Item($row, 6).Value = 133553
$ws.Cells.Item($row, 7).Value = 158433
$ws.Cells.Item($row, 8).Value = 130842
$ws.Cells.Item($row, 9).Value = 174624
$ws.Cells.Item($row, 10).Value = 113750
$ws.Cells.Item($row, 11).Value = 198644
$ws.Cells.Item($row, 12).Value = 220588
$ws.Cells.Item($row, 13).Value = 171963
$ws.Cells.Item($row, 14).Value = 119871
$ws.Cells.Item($row, 15).Value = 38707
$ws.Cells.Item($row, 16).Value = 34826
$ws.Cells.Item($row, 17).Value = 50540
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36931
$ws.Cells.Item($row, 20).Value = -1
